$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Poroto verde" (Feria Lagunitas
# de Puerto Montt) series. Insert a fresh row above the existing row 19 so
# every subsequent record shifts down by one (old row 19 -> new row 20,
# ..., old row 79 -> new row 80), then populate the newly inserted row 19
# with the new record's data.
$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44708
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = 100112031
$ws.Cells.Item(19, 7).Value = "Poroto verde"
$ws.Cells.Item(19, 8).Value = "Magnum"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 25000
$ws.Cells.Item(19, 12).Value = 25000
$ws.Cells.Item(19, 13).Value = 25000
$ws.Cells.Item(19, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 1000
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
